$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(1,1).Range.Text = "54-29=25"
$t.Cell(1,2).Range.Text = "12-9=3"
$t.Cell(1,3).Range.Text = "30+25=55"
$t.Cell(1,4).Range.Text = "85+14=99"
$t.Cell(1,5).Range.Text = "83-67=16"
$t.Cell(2,1).Range.Text = "42+36=78"
$t.Cell(2,2).Range.Text = "90-42=48"
$t.Cell(2,3).Range.Text = "47-43=4"
$t.Cell(2,4).Range.Text = "72-54=18"
$t.Cell(2,5).Range.Text = "72-8=64"
$t.Cell(3,1).Range.Text = "45+48=93"
$t.Cell(3,2).Range.Text = "52+7=59"
$t.Cell(3,3).Range.Text = "9+21=30"
$t.Cell(3,4).Range.Text = "38+30=68"
$t.Cell(3,5).Range.Text = "8+32=40"
$t.Cell(4,1).Range.Text = "8+34=42"
$t.Cell(4,2).Range.Text = "83+16=99"
$t.Cell(4,3).Range.Text = "31-20=11"
$t.Cell(4,4).Range.Text = "27-5=22"
$t.Cell(4,5).Range.Text = "17+79=96"
$t.Cell(5,1).Range.Text = "95-3=92"
$t.Cell(5,2).Range.Text = "84+11=95"
$t.Cell(5,3).Range.Text = "75-16=59"
$t.Cell(5,4).Range.Text = "70+23=93"
$t.Cell(5,5).Range.Text = "17+78=95"
$t.Cell(6,1).Range.Text = "28+38=66"
$t.Cell(6,2).Range.Text = "50+35=85"
$t.Cell(6,3).Range.Text = "5+65=70"
$t.Cell(6,4).Range.Text = "24+54=78"
$t.Cell(6,5).Range.Text = "44+7=51"
$t.Cell(7,1).Range.Text = "98-71=27"
$t.Cell(7,2).Range.Text = "25+34=59"
$t.Cell(7,3).Range.Text = "92-1=91"
$t.Cell(7,4).Range.Text = "94-86=8"
$t.Cell(7,5).Range.Text = "89-2=87"
$t.Cell(8,1).Range.Text = "62+17=79"
$t.Cell(8,2).Range.Text = "82-35=47"
$t.Cell(8,3).Range.Text = "86+1=87"
$t.Cell(8,4).Range.Text = "56+17=73"
$t.Cell(8,5).Range.Text = "97-67=30"
$t.Cell(9,1).Range.Text = "42-39=3"
$t.Cell(9,2).Range.Text = "20+50=70"
$t.Cell(9,3).Range.Text = "32-12=20"
$t.Cell(9,4).Range.Text = "71-7=64"
$t.Cell(9,5).Range.Text = "33+17=50"
$t.Cell(10,1).Range.Text = "9+66=75"
$t.Cell(10,2).Range.Text = "42+36=78"
$t.Cell(10,3).Range.Text = "26+63=89"
$t.Cell(10,4).Range.Text = "2+92=94"
$t.Cell(10,5).Range.Text = "10+69=79"
$t.Cell(11,1).Range.Text = "47+1=48"
$t.Cell(11,2).Range.Text = "33+60=93"
$t.Cell(11,3).Range.Text = "30+34=64"
$t.Cell(11,4).Range.Text = "0+28=28"
$t.Cell(11,5).Range.Text = "52+46=98"
$t.Cell(12,1).Range.Text = "55+6=61"
$t.Cell(12,2).Range.Text = "15+61=76"
$t.Cell(12,3).Range.Text = "83-8=75"
$t.Cell(12,4).Range.Text = "6+11=17"
$t.Cell(12,5).Range.Text = "33-12=21"
$t.Cell(13,1).Range.Text = "7+9=16"
$t.Cell(13,2).Range.Text = "79-64=15"
$t.Cell(13,3).Range.Text = "13+74=87"
$t.Cell(13,4).Range.Text = "95-86=9"
$t.Cell(13,5).Range.Text = "57-38=19"
$t.Cell(14,1).Range.Text = "89-1=88"
$t.Cell(14,2).Range.Text = "47+51=98"
$t.Cell(14,3).Range.Text = "14+76=90"
$t.Cell(14,4).Range.Text = "69-69=0"
$t.Cell(14,5).Range.Text = "27+46=73"
$t.Cell(15,1).Range.Text = "99-86=13"
$t.Cell(15,2).Range.Text = "3+68=71"
$t.Cell(15,3).Range.Text = "59+15=74"
$t.Cell(15,4).Range.Text = "63+27=90"
$t.Cell(15,5).Range.Text = "51+0=51"
$t.Cell(16,1).Range.Text = "24+37=61"
$t.Cell(16,2).Range.Text = "63+8=71"
$t.Cell(16,3).Range.Text = "68-21=47"
$t.Cell(16,4).Range.Text = "11-0=11"
$t.Cell(16,5).Range.Text = "74-48=26"
$t.Cell(17,1).Range.Text = "77-10=67"
$t.Cell(17,2).Range.Text = "67+14=81"
$t.Cell(17,3).Range.Text = "24+63=87"
$t.Cell(17,4).Range.Text = "99-93=6"
$t.Cell(17,5).Range.Text = "56+2=58"
$t.Cell(18,1).Range.Text = "7+45=52"
$t.Cell(18,2).Range.Text = "41-15=26"
$t.Cell(18,3).Range.Text = "12+50=62"
$t.Cell(18,4).Range.Text = "87-43=44"
$t.Cell(18,5).Range.Text = "22+67=89"
$t.Cell(19,1).Range.Text = "35+49=84"
$t.Cell(19,2).Range.Text = "35+51=86"
$t.Cell(19,3).Range.Text = "70-69=1"
$t.Cell(19,4).Range.Text = "50-39=11"
$t.Cell(19,5).Range.Text = "74-37=37"
$t.Cell(20,1).Range.Text = "36+37=73"
$t.Cell(20,2).Range.Text = "3+35=38"
$t.Cell(20,3).Range.Text = "77-42=35"
$t.Cell(20,4).Range.Text = "61-28=33"
$t.Cell(20,5).Range.Text = "67-41=26"
